# Update countries & provincias Spain
#
# Applies the data refresh to the "Pais" sheet of paises.xlsx:
#  - Argentina overtakes Kuwait (rows 56-57) with new figures
#  - "Consejo Danes para los Refugiados" overtakes Somalia, Uruguay,
#    Burkina Faso and Guatemala (rows 104-108), each shifting down one row
#  - Bahamas overtakes Guyana (rows 159-160) with new figures
#  - Assorted numeric updates (Estados Unidos totals row, El Salvador,
#    Monaco)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Pais")

function Set-CountryRow($r, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $country
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Row 4: Estados Unidos - refreshed totals
Set-CountryRow 4 "Estados Unidos" 1159923 28893 160705 931818 16475 1647 67400

# Rows 56-57: Argentina overtakes Kuwait
Set-CountryRow 56 "Argentina" 4681 149 1320 3124 157 12 237
Set-CountryRow 57 "Kuwait" 4619 242 1703 2883 69 3 33

# Rows 104-108: "Consejo Danes para los Refugiados" overtakes Somalia,
# Uruguay, Burkina Faso and Guatemala, each shifting down one row
Set-CountryRow 104 "Consejo Danes para los Refugiados" 674 70 75 566 0 1 33
Set-CountryRow 105 "Somalia" 671 70 34 606 2 3 31
Set-CountryRow 106 "Uruguay" 652 4 440 195 10 0 17
Set-CountryRow 107 "Burkina Faso" 652 3 535 73 0 0 44
Set-CountryRow 108 "Guatemala" 644 45 72 556 5 0 16

# Row 118: El Salvador - refreshed figures
$ws.Cells.Item(118, 5).Value = 294
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 11

# Row 155: Monaco - refreshed figures
$ws.Cells.Item(155, 4).Value = 78
$ws.Cells.Item(155, 5).Value = 13

# Rows 159-160: Bahamas overtakes Guyana
Set-CountryRow 159 "Bahamas" 83 2 24 48 1 0 11
Set-CountryRow 160 "Guyana" 82 0 22 51 2 0 9
